# stemming_rules.xlsx update
# - Reworks the "ata" stemming rule (row 33) into a regex-based rule with a
#   back-reference replacement and an explanatory comment.
# - Adds four new stemming rules at the bottom of the table (rows 69-72):
#   wala, khani, ma?nda?, ika?.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: "ata" rule -> regex-based rule ------------------------------
$ws.Range("A33").Value = "(^\w.*?[^aeiou][ai]?)ta"
$ws.Range("B33").Value = "\1"
$ws.Range("C33").Value = "Adhunikata , adhunik  [Change only if it is preceded by 2 consonants(ie: syllables).... eg (bata !-> ba ) X]  (visheshta -> vishesh) , (pravahita -> pravah)"

# --- New rules appended at the bottom of the table ------------------------
$ws.Range("A69").Value = "wala"
$ws.Range("B69").Value = "_"
$ws.Range("C69").Value = "bichwala -> bich"

$ws.Range("A70").Value = "khani"
$ws.Range("B70").Value = "_"
$ws.Range("C70").Value = "Chedkhani -> cheda"

$ws.Range("A71").Value = "ma?nda?"
$ws.Range("B71").Value = "_"
$ws.Range("C71").Value = "Zaruratmand -> zarurat"

$ws.Range("A72").Value = "ika?"
$ws.Range("B72").Value = "_"
$ws.Range("C72").Value = "Itihasik -> itihas"

# The original edit left cell A70 with direct formatting applied (same Arial
# font, just marked as explicitly-applied) - reapply the cell style so the
# same xf bookkeeping is produced.
$ws.Range("A70").Style = "Normal"

# --- Move the selection to the last edited cell, like the author did ------
$ws.Range("A72").Select()
